$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.021.77'
$ws.Range('E2').Value = '  +0.75%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.904.01'
$ws.Range('E3').Value = '  +0.39%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8352'
$ws.Range('E5').Value = '  +9.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.89'
$ws.Range('E6').Value = '  +0.74%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3232'
$ws.Range('E8').Value = '  +6.51%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.96'
$ws.Range('E9').Value = '  +6.38%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07028'
$ws.Range('E10').Value = '  +3.08%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08027'
$ws.Range('E11').Value = '  +0.70%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7522'
$ws.Range('E12').Value = '  +2.51%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.887.05'
$ws.Range('E13').Value = '  -0.04%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.221'
$ws.Range('E14').Value = '  +1.36%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.80'
$ws.Range('E15').Value = '  +2.17%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.034.84'
$ws.Range('E16').Value = '  +0.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.17'
$ws.Range('E17').Value = '  +2.84%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.947'
$ws.Range('E18').Value = '  +0.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.66'
$ws.Range('E19').Value = '  +1.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007783'
$ws.Range('E20').Value = '  +1.24%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.164.03'
$ws.Range('E21').Value = '  +1.43%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.997'
$ws.Range('E24').Value = '  +1.86%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1594'
$ws.Range('E25').Value = '  +24.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '168.59'
$ws.Range('E26').Value = '  +1.23%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.231'
$ws.Range('E27').Value = '  +0.30%  '

$ws.Range('E28').Value = '  +2.04%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.094'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.372'
$ws.Range('E30').Value = '  -1.96%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.522'
$ws.Range('E31').Value = '  +0.66%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.285'
$ws.Range('E32').Value = '  +0.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05750'
$ws.Range('E33').Value = '  +10.74%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.094'
$ws.Range('E34').Value = '  +0.89%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.293'
$ws.Range('E35').Value = '  +3.94%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7359'
$ws.Range('E36').Value = '  +1.81%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.725'
$ws.Range('E37').Value = '  +0.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01910'
$ws.Range('E38').Value = '  -0.06%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.795'
$ws.Range('E39').Value = '  +0.83%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4430'
$ws.Range('E40').Value = '  +1.10%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.33'
$ws.Range('E41').Value = '  +0.90%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.963'
$ws.Range('E42').Value = '  -2.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8447'
$ws.Range('E43').Value = '  +2.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('E44').Value = '  -0.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.901'
$ws.Range('E45').Value = '  +1.15%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.634'
$ws.Range('E46').Value = '  +0.60%  '

$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.40'
$ws.Range('E47').Value = '  +1.90%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.734'
$ws.Range('E48').Value = '  -0.19%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '991.57'
$ws.Range('E49').Value = '  +9.20%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.058.81'
$ws.Range('E50').Value = '  +1.22%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.28'
$ws.Range('E51').Value = '  +0.67%  '
